$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1353.6666
$ws.Range("J2").Value = 5749.5
$ws.Range("L2").Value = 5749.5
$ws.Range("N2").Value = -5975.5
$ws.Range("H74").Value = 10578.044
$ws.Range("I74").Value = 8858.083000000001
$ws.Range("J74").Value = 12454.363
$ws.Range("K74").Value = 8858.083000000001
$ws.Range("L74").Value = 12454.363
$ws.Range("M74").Value = -7922.083000000001
$ws.Range("N74").Value = -14326.363
$ws.Range("H77").Value = 10578.044
$ws.Range("I77").Value = 8858.083000000001
$ws.Range("J77").Value = 12454.363
$ws.Range("K77").Value = 44290.415
$ws.Range("L77").Value = 62271.815
$ws.Range("M77").Value = -39610.415
$ws.Range("N77").Value = -71631.815
$ws.Range("H87").Value = 60000
$ws.Range("J87").Value = 60000
$ws.Range("L87").Value = 60000
$ws.Range("N87").Value = -62496
$ws.Range("H90").Value = 60000
$ws.Range("J90").Value = 60000
$ws.Range("L90").Value = 180000
$ws.Range("N90").Value = -192480
$ws.Range("H132").Value = 1882.1892
$ws.Range("I132").Value = 1420.8387
$ws.Range("K132").Value = 4262.5161
$ws.Range("M132").Value = -1732.5161
$ws.Range("H137").Value = 4779.5415
$ws.Range("I137").Value = 5240.45
$ws.Range("K137").Value = 15721.35
$ws.Range("M137").Value = -13171.35

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7922.185
$ws.Range("I61").Value = 8034.577
$ws.Range("K61").Value = 8034.577
$ws.Range("M61").Value = -7822.577
$ws.Range("H122").Value = 2666.6667
$ws.Range("J122").Value = 3000
$ws.Range("L122").Value = 9000
$ws.Range("N122").Value = -13900
$ws.Range("H132").Value = 2347.7827
$ws.Range("I132").Value = 2181.2
$ws.Range("K132").Value = 6543.599999999999
$ws.Range("M132").Value = -4013.599999999999
$ws.Range("H136").Value = 7922.185
$ws.Range("I136").Value = 8034.577
$ws.Range("K136").Value = 24103.731
$ws.Range("M136").Value = -21553.731

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3700
$ws.Range("I86").Value = 3700
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 3700
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = ""
$ws.Range("N86").Value = -2577
$ws.Range("H89").Value = 3700
$ws.Range("I89").Value = 3700
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 18500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = ""
$ws.Range("N89").Value = -12884
$ws.Range("H94").Value = 671.5
$ws.Range("I94").Value = 556.0625
$ws.Range("J94").Value = 1595
$ws.Range("K94").Value = 556.0625
$ws.Range("L94").Value = 1595
$ws.Range("M94").Value = -105.0625
$ws.Range("N94").Value = -2497
$ws.Range("H107").Value = 1301.5518
$ws.Range("I107").Value = 1342
$ws.Range("J107").Value = 1195.375
$ws.Range("K107").Value = 1342
$ws.Range("L107").Value = 1195.375
$ws.Range("M107").Value = 578
$ws.Range("N107").Value = -5035.375
$ws.Range("H134").Value = 7007.0728
$ws.Range("I134").Value = 6816.265
$ws.Range("J134").Value = 8565.333000000001
$ws.Range("K134").Value = 20448.795
$ws.Range("L134").Value = 25695.999
$ws.Range("M134").Value = -17913.795
$ws.Range("N134").Value = -30765.999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 13861.909
$ws.Range("J59").Value = 11248.1
$ws.Range("L59").Value = 11248.1
$ws.Range("N59").Value = -13538.1
$ws.Range("H74").Value = 40249.668
$ws.Range("J74").Value = 40249.668
$ws.Range("L74").Value = 40249.668
$ws.Range("N74").Value = -41997.668
$ws.Range("H77").Value = 40249.668
$ws.Range("J77").Value = 40249.668
$ws.Range("L77").Value = 120749.004
$ws.Range("N77").Value = -129485.004
$ws.Range("H105").Value = 2166.3333
$ws.Range("I105").Value = 2250
$ws.Range("K105").Value = 2250
$ws.Range("M105").Value = -503
$ws.Range("H122").Value = 3250.2222
$ws.Range("I122").Value = 3427.3333
$ws.Range("J122").Value = 1833.3334
$ws.Range("K122").Value = 10281.9999
$ws.Range("L122").Value = 5500.0002
$ws.Range("M122").Value = -7831.999899999999
$ws.Range("N122").Value = -10400.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 140725.47
$ws.Range("J37").Value = 140725.47
$ws.Range("L37").Value = 422176.41
$ws.Range("N37").Value = -422400.41
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = ""
$ws.Range("M41").Value = ""
$ws.Range("N41").Value = 0
$ws.Range("H68").Value = 45459864
$ws.Range("I68").Value = 125002620
$ws.Range("K68").Value = 375007860
$ws.Range("M68").Value = -375007049
$ws.Range("H71").Value = 45459864
$ws.Range("I71").Value = 125002620
$ws.Range("K71").Value = 1125023580
$ws.Range("M71").Value = -1125019524
$ws.Range("H107").Value = 328.125
$ws.Range("I107").Value = 295
$ws.Range("J107").Value = 361.25
$ws.Range("K107").Value = 885
$ws.Range("L107").Value = 1083.75
$ws.Range("M107").Value = 1035
$ws.Range("N107").Value = -4923.75
$ws.Range("H109").Value = 13709.1
$ws.Range("I109").Value = 30684.666
$ws.Range("K109").Value = 92053.99800000001
$ws.Range("M109").Value = -91013.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H99").Value = 8914.9
$ws.Range("J99").Value = 40000
$ws.Range("L99").Value = 40000
$ws.Range("N99").Value = -44492
$ws.Range("I122").Value = 3424.25
$ws.Range("J122").Value = 2125
$ws.Range("K122").Value = 10272.75
$ws.Range("L122").Value = 6375
$ws.Range("M122").Value = -7822.75
$ws.Range("N122").Value = -11275
$ws.Range("H126").Value = 6186.5
$ws.Range("I126").Value = 4624.5
$ws.Range("J126").Value = 7748.5
$ws.Range("K126").Value = 13873.5
$ws.Range("L126").Value = 23245.5
$ws.Range("M126").Value = -11403.5
$ws.Range("N126").Value = -28185.5
$ws.Range("H132").Value = 4354.606
$ws.Range("I132").Value = 3789.1667
$ws.Range("K132").Value = 11367.5001
$ws.Range("M132").Value = -8837.500100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 297.7857
$ws.Range("I16").Value = 297.7857
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 297.7857
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = ""
$ws.Range("N16").Value = -127.7857
$ws.Range("H46").Value = 6807.8623
$ws.Range("J46").Value = 7099.875
$ws.Range("L46").Value = 7099.875
$ws.Range("N46").Value = -7475.875
$ws.Range("H68").Value = 10248
$ws.Range("I68").Value = 9996.5
$ws.Range("K68").Value = 9996.5
$ws.Range("M68").Value = -9247.5
$ws.Range("H71").Value = 10248
$ws.Range("I71").Value = 9996.5
$ws.Range("K71").Value = 49982.5
$ws.Range("M71").Value = -46238.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1602.1666
$ws.Range("I81").Value = 1475.091
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 2950.182
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -1889.182
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 1602.1666
$ws.Range("I84").Value = 1475.091
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 14750.91
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -9446.91
$ws.Range("N84").Value = -40608
$ws.Range("H107").Value = 1034.3
$ws.Range("I107").Value = 1210.2307
$ws.Range("J107").Value = 707.5714
$ws.Range("K107").Value = 3630.6921
$ws.Range("L107").Value = 2122.7142
$ws.Range("M107").Value = -1710.6921
$ws.Range("N107").Value = -5962.7142
$ws.Range("H113").Value = 1651.9
$ws.Range("I113").Value = 702.8333
$ws.Range("K113").Value = 2108.4999
$ws.Range("M113").Value = 61.5001000000002
$ws.Range("H122").Value = 2787.9697
$ws.Range("I122").Value = 1888.48
$ws.Range("J122").Value = 5598.875
$ws.Range("K122").Value = 5665.440000000001
$ws.Range("L122").Value = 16796.625
$ws.Range("M122").Value = -3215.440000000001
$ws.Range("N122").Value = -21696.625
$ws.Range("H136").Value = 7665.778
$ws.Range("J136").Value = 10247.25
$ws.Range("L136").Value = 30741.75
$ws.Range("N136").Value = -35841.75
